$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.568.64"
$ws.Range("E2").Value = "  -2.66%  "

# Row 3
$ws.Range("D3").Value = "2.004.11"
$ws.Range("E3").Value = "  -4.21%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +1.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.87"
$ws.Range("E5").Value = "  -3.75%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.012"
$ws.Range("E6").Value = "  +1.11%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5006"
$ws.Range("E7").Value = "  -4.36%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4221"
$ws.Range("E8").Value = "  -4.60%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.64"
$ws.Range("E9").Value = "  +0.29%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09028"
$ws.Range("E10").Value = "  -3.28%  "

# Row 11
$ws.Range("E11").Value = "  -4.39%  "

# Row 12
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.29"
$ws.Range("E12").Value = "  -6.32%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "2.074.40"
$ws.Range("E13").Value = "  +2.13%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.046"
$ws.Range("E14").Value = "  -6.45%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.466"
$ws.Range("E15").Value = "  -6.29%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.015"
$ws.Range("E16").Value = "  +1.35%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.43"
$ws.Range("E17").Value = "  -6.78%  "

# Row 18
$ws.Range("E18").Value = "  -3.76%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06680"
$ws.Range("E19").Value = "  +0.24%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.65"
$ws.Range("E20").Value = "  -7.07%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.012"
$ws.Range("E21").Value = "  +1.18%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.970"
$ws.Range("E22").Value = "  -5.71%  "

# Row 23
$ws.Range("D23").Value = "29.607.16"
$ws.Range("E23").Value = "  -2.43%  "

# Row 24
$ws.Range("E24").Value = "  -4.46%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.306"
$ws.Range("E25").Value = "  +0.27%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.75"
$ws.Range("E26").Value = "  -2.42%  "

# Row 27
$ws.Range("E27").Value = "  -4.95%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.346"
$ws.Range("E28").Value = "  -6.14%  "

# Row 29
$ws.Range("E29").Value = "  -8.66%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.18"
$ws.Range("E30").Value = "  -3.67%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.057"
$ws.Range("E31").Value = "  -7.38%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09951"
$ws.Range("E32").Value = "  -4.82%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.563"
$ws.Range("E33").Value = "  -6.71%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.838"
$ws.Range("E34").Value = "  -6.57%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.802"

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02465"
$ws.Range("E36").Value = "  -6.42%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.267"

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.307"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06412"
$ws.Range("E39").Value = "  -6.22%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6552"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.68"
$ws.Range("E41").Value = "  -6.70%  "

# Row 42
$ws.Range("E42").Value = "  -7.36%  "

# Row 43
$ws.Range("E43").Value = "  +1.20%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6348"
$ws.Range("E44").Value = "  -6.99%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.44"
$ws.Range("E45").Value = "  -6.74%  "

# Row 46
$ws.Range("E46").Value = "  -6.60%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.301"
$ws.Range("E47").Value = "  -5.15%  "

# Row 48
$ws.Range("E48").Value = "  -3.34%  "

# Row 49
$ws.Range("E49").Value = "  -1.35%  "

# Row 50
$ws.Range("E50").Value = "  -3.44%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.128"
$ws.Range("E51").Value = "  -6.75%  "
